$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H80").Value = 1821.2858
$ws.Range("I80").Value = 1700
$ws.Range("J80").Value = 1983
$ws.Range("K80").Value = 5100
$ws.Range("L80").Value = 5949
$ws.Range("M80").Value = -4102
$ws.Range("N80").Value = -7945
$ws.Range("H83").Value = 1821.2858
$ws.Range("I83").Value = 1700
$ws.Range("J83").Value = 1983
$ws.Range("K83").Value = 15300
$ws.Range("L83").Value = 17847
$ws.Range("M83").Value = -10308
$ws.Range("N83").Value = -27831
$ws.Range("H137").Value = 4192.6
$ws.Range("J137").Value = 8500
$ws.Range("L137").Value = 25500
$ws.Range("N137").Value = -30600
$ws.Range("H138").Value = 2155.923
$ws.Range("I138").Value = 1432.4286
$ws.Range("K138").Value = 4297.2858
$ws.Range("M138").Value = 842.7142000000003

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()
$ws.Range("H61").Value = 2812
$ws.Range("I61").Value = 2812
$ws.Range("J61").Value = 0
$ws.Range("K61").Value = 2812
$ws.Range("L61").Value = 0
$ws.Range("M61").Value = -2600
$ws.Range("N61").ClearContents()
$ws.Range("H116").Value = 0
$ws.Range("I116").Value = 0
$ws.Range("K116").Value = 0
$ws.Range("M116").ClearContents()
$ws.Range("H136").Value = 2812
$ws.Range("I136").Value = 2812
$ws.Range("J136").Value = 0
$ws.Range("K136").Value = 8436
$ws.Range("L136").Value = 0
$ws.Range("M136").Value = -5886
$ws.Range("N136").ClearContents()

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 0
$ws.Range("I3").Value = 0
$ws.Range("K3").Value = 0
$ws.Range("M3").ClearContents()
$ws.Range("H86").Value = 4813.875
$ws.Range("I86").Value = 4644.4287
$ws.Range("J86").Value = 6000
$ws.Range("K86").Value = 4644.4287
$ws.Range("L86").Value = 6000
$ws.Range("M86").Value = -3521.4287
$ws.Range("N86").Value = -8246
$ws.Range("H89").Value = 4813.875
$ws.Range("I89").Value = 4644.4287
$ws.Range("J89").Value = 6000
$ws.Range("K89").Value = 23222.1435
$ws.Range("L89").Value = 30000
$ws.Range("M89").Value = -17606.1435
$ws.Range("N89").Value = -41232

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2897.111
$ws.Range("I16").Value = 2358.3333
$ws.Range("J16").Value = 3974.6667
$ws.Range("K16").Value = 2358.3333
$ws.Range("L16").Value = 3974.6667
$ws.Range("M16").Value = -2071.3333
$ws.Range("N16").Value = -4548.6667
$ws.Range("H113").Value = 2897.111
$ws.Range("I113").Value = 2358.3333
$ws.Range("J113").Value = 3974.6667
$ws.Range("K113").Value = 2358.3333
$ws.Range("L113").Value = 3974.6667
$ws.Range("M113").Value = -188.3332999999998
$ws.Range("N113").Value = -8314.6667
$ws.Range("H121").Value = 0
$ws.Range("J121").Value = 0
$ws.Range("L121").Value = 0
$ws.Range("N121").ClearContents()
$ws.Range("H132").Value = 2253.125
$ws.Range("I132").Value = 2100.3845
$ws.Range("J132").Value = 2915
$ws.Range("K132").Value = 6301.1535
$ws.Range("L132").Value = 8745
$ws.Range("M132").Value = -3771.1535
$ws.Range("N132").Value = -13805

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H2").Value = 33.5
$ws.Range("I2").Value = 40.125
$ws.Range("J2").Value = 24.666666
$ws.Range("K2").Value = 240.75
$ws.Range("L2").Value = 147.999996
$ws.Range("M2").Value = -127.75
$ws.Range("N2").Value = -373.999996
$ws.Range("H34").Value = 965.38464
$ws.Range("I34").Value = 520
$ws.Range("J34").Value = 1243.75
$ws.Range("K34").Value = 1560
$ws.Range("L34").Value = 3731.25
$ws.Range("M34").Value = -1476
$ws.Range("N34").Value = -3899.25
$ws.Range("H38").Value = 1141
$ws.Range("I38").Value = 1917.25
$ws.Range("J38").Value = 106
$ws.Range("K38").Value = 5751.75
$ws.Range("L38").Value = 318
$ws.Range("M38").Value = -5404.75
$ws.Range("N38").Value = -1012
$ws.Range("H40").Value = 11
$ws.Range("I40").Value = 15
$ws.Range("J40").Value = 5
$ws.Range("K40").Value = 60
$ws.Range("L40").Value = 20
$ws.Range("M40").Value = 9
$ws.Range("N40").Value = -158

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H68").Value = 659.6
$ws.Range("I68").Value = 659.6
$ws.Range("K68").Value = 659.6
$ws.Range("M68").Value = 89.39999999999998
$ws.Range("H71").Value = 659.6
$ws.Range("I71").Value = 659.6
$ws.Range("K71").Value = 3298
$ws.Range("M71").Value = 446
$ws.Range("H82").Value = 1953.5
$ws.Range("I82").Value = 2499
$ws.Range("K82").Value = 2499
$ws.Range("M82").Value = -2138
$ws.Range("H85").Value = 1953.5
$ws.Range("I85").Value = 2499
$ws.Range("K85").Value = 2499
$ws.Range("M85").Value = -1251

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H62").Value = 4791.273
$ws.Range("I62").Value = 3496.6667
$ws.Range("J62").Value = 5276.75
$ws.Range("K62").Value = 3496.6667
$ws.Range("L62").Value = 5276.75
$ws.Range("M62").Value = -2872.6667
$ws.Range("N62").Value = -6524.75
$ws.Range("H65").Value = 4791.273
$ws.Range("I65").Value = 3496.6667
$ws.Range("J65").Value = 5276.75
$ws.Range("K65").Value = 17483.3335
$ws.Range("L65").Value = 26383.75
$ws.Range("M65").Value = -14363.3335
$ws.Range("N65").Value = -32623.75
$ws.Range("H82").Value = 30150.5
$ws.Range("J82").Value = 30150.5
$ws.Range("L82").Value = 30150.5
$ws.Range("N82").Value = -30916.5
$ws.Range("H85").Value = 30150.5
$ws.Range("J85").Value = 30150.5
$ws.Range("L85").Value = 30150.5
$ws.Range("N85").Value = -32802.5
$ws.Range("H132").Value = 1654.1666
$ws.Range("I132").Value = 1632.125
$ws.Range("J132").Value = 1698.25
$ws.Range("K132").Value = 4896.375
$ws.Range("L132").Value = 5094.75
$ws.Range("M132").Value = -2366.375
$ws.Range("N132").Value = -10154.75
